$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("L2:L45").Value = 3.0
